# Add a 7th dataset entry ("bank-marketing-dataset-analysis-classification")
# to the "6 лаба" row: put the dataset name in C8 as a hyperlink (matching
# the style already used for the other dataset links), grow row 8 to the
# same height as the other hyperlinked rows, and leave the selection on C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = "bank-marketing-dataset-analysis-classification"

$ws.Hyperlinks.Add(
    $ws.Range("C8"),
    "https://www.kaggle.com/datasets/janiobachmann/bank-marketing-dataset"
)

# Hyperlinks.Add re-derives a fresh cell style; restore the shared
# "hyperlink" formatting already used by the sibling cells (e.g. C7) so
# C8 keeps using the same style as the rest of column C.
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(8).RowHeight = 45

$ws.Range("C8").Select()
